# Progress-journal edit:
#   "***need to create getAcceptedTrade…. Functions****" (with the
#   _GoBack bookmark) is turned into a new "running out of gas" note,
#   and the original "need to create getAcceptedTrade" note is kept,
#   re-appearing as its own paragraph right after (without the bookmark).

$d = $word.ActiveDocument

# --- locate the target paragraph dynamically (avoid a brittle hard-coded index) ---
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*need to create*getAcceptedTrade*Functions*") {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq -1) {
    throw "Could not locate the 'need to create getAcceptedTrade' paragraph"
}

$p = $d.Paragraphs.Item($targetIndex)

# --- duplicate the paragraph's content into a new paragraph right after it ---
# Assigning FormattedText at the paragraph's own end (just before its
# paragraph mark) makes Word append a paragraph break + the duplicated runs
# immediately after, instead of leaving a stray empty run in a pre-inserted
# empty paragraph.
$origEnd = $p.Range.End
$dup = $p.Range.FormattedText
$pasteTarget = $d.Range($origEnd, $origEnd)
$pasteTarget.FormattedText = $dup

# The bookmark travelled with the original (first) paragraph, which is what
# we want: the duplicate (second, new) paragraph should NOT carry it.

# --- rewrite the text of the original (now first) paragraph in place ---
# Editing the existing runs' text via Find/Replace (rather than replacing the
# whole paragraph range) preserves the surrounding run/proofErr structure.
$editRange = $d.Paragraphs.Item($targetIndex).Range
$editRange.Find.Execute("***need to create ", $false, $false, $false, $false, $false, $true, 0, $false, "***running into problem: running out of gas when calling ", 2)

$editRange = $d.Paragraphs.Item($targetIndex).Range
$editRange.Find.Execute("getAcceptedTrade", $false, $false, $false, $false, $false, $true, 0, $false, "searchByStartDate", 2)

$editRange = $d.Paragraphs.Item($targetIndex).Range
$editRange.Find.Execute([char]0x2026 + ". Functions****", $false, $false, $false, $false, $false, $true, 0, $false, " function***", 2)
